$d = $word.ActiveDocument

# The document's title and abstract paragraphs were originally split into
# many single-word/space runs (e.g. one <w:r> per word). The edit merges
# each of those paragraphs' text into a single run by replacing the full
# paragraph text in place (Find/Replace coalesces the matched range into
# one run holding the whole replacement string).

$d.Content.Find.Execute(
    "Questions: Trigonometry (degrees)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Questions: Trigonometry (degrees)", 2
)

$d.Content.Find.Execute(
    "A selection of questions on trigonometry, where angles are measured in degrees.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A selection of questions on trigonometry, where angles are measured in degrees.", 2
)
